$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("chart")

$data = @(
    @(3819, 3830, 318, 2190, 2093, 106),
    @(4777, 4752, 318, 2154, 2061, 108),
    @(4158, 3823, 590, 3081, 2537, 599),
    @(3969, 3825, 428, 4123, 4044, 119),
    @(4909, 4748, 431, 4110, 4006, 118),
    @(4092, 3796, 588, 5234, 4900, 420),
    @(3925, 3741, 456, 4178, 4033, 128),
    @(4015, 3830, 452, 4196, 4078, 142),
    @(4697, 3893, 1106, 5622, 4763, 928),
    @(5180, 4741, 721, 4165, 4039, 178),
    @(4996, 4743, 549, 4170, 4010, 150),
    @(4113, 3735, 635, 4211, 4071, 185),
    @(5019, 4757, 551, 8204, 7939, 290),
    @(5122, 4754, 642, 8411, 8207, 169)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
